$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$values = @{
    2  = "2025-10-30T03:38:27.802493+00:00"
    3  = "2025-10-30T03:38:29.736682+00:00"
    4  = "2025-10-30T03:38:29.736698+00:00"
    5  = "2025-10-30T03:38:29.736707+00:00"
    6  = "2025-10-30T03:38:29.736715+00:00"
    7  = "2025-10-30T03:38:29.736722+00:00"
    8  = "2025-10-30T03:38:31.706255+00:00"
    9  = "2025-10-30T03:38:31.706271+00:00"
    10 = "2025-10-30T03:38:33.678259+00:00"
    11 = "2025-10-30T03:38:35.653629+00:00"
    12 = "2025-10-30T03:38:35.653645+00:00"
    13 = "2025-10-30T03:38:38.110253+00:00"
    14 = "2025-10-30T03:38:38.110269+00:00"
    15 = "2025-10-30T03:38:38.110276+00:00"
    16 = "2025-10-30T03:38:38.110284+00:00"
    17 = "2025-10-30T03:38:45.519642+00:00"
    18 = "2025-10-30T03:38:47.573574+00:00"
    19 = "2025-10-30T03:38:49.522526+00:00"
    20 = "2025-10-30T03:38:51.937235+00:00"
    21 = "2025-10-30T03:38:51.937251+00:00"
    22 = "2025-10-30T03:38:51.937259+00:00"
    23 = "2025-10-30T03:38:54.307093+00:00"
    24 = "2025-10-30T03:38:54.307109+00:00"
    25 = "2025-10-30T03:38:54.307116+00:00"
    26 = "2025-10-30T03:38:54.307124+00:00"
    27 = "2025-10-30T03:38:54.307131+00:00"
    28 = "2025-10-30T03:39:01.507967+00:00"
    29 = "2025-10-30T03:39:01.507984+00:00"
    30 = "2025-10-30T03:39:01.507993+00:00"
    31 = "2025-10-30T03:39:01.508001+00:00"
    32 = "2025-10-30T03:39:03.978508+00:00"
    33 = "2025-10-30T03:39:03.978525+00:00"
    34 = "2025-10-30T03:39:03.978532+00:00"
    35 = "2025-10-30T03:39:06.146788+00:00"
    36 = "2025-10-30T03:39:06.146806+00:00"
    37 = "2025-10-30T03:39:06.146813+00:00"
    38 = "2025-10-30T03:39:06.146820+00:00"
    39 = "2025-10-30T03:39:06.146827+00:00"
    40 = "2025-10-30T03:39:06.146833+00:00"
    41 = "2025-10-30T03:39:06.146840+00:00"
    42 = "2025-10-30T03:39:06.146850+00:00"
    43 = "2025-10-30T03:39:06.146856+00:00"
    44 = "2025-10-30T03:39:08.154979+00:00"
    45 = "2025-10-30T03:39:08.154995+00:00"
    46 = "2025-10-30T03:39:12.579182+00:00"
    47 = "2025-10-30T03:39:14.567908+00:00"
    48 = "2025-10-30T03:39:14.567923+00:00"
    49 = "2025-10-30T03:39:14.567931+00:00"
    50 = "2025-10-30T03:39:14.567937+00:00"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 11).Value = $values[$row]
}
